$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = 0.827942324055343
$ws.Range("C6").Value = 0.02598578151053643
$ws.Range("D6").Value = 0.7898560923698297
$ws.Range("E6").Value = 0.8200384473729125
$ws.Range("F6").Value = 0.8294487770294222
$ws.Range("G6").Value = 0.8293082722044585
$ws.Range("H6").Value = 0.8710600313000915
$ws.Range("I6").Value = 0.8142274241189149
$ws.Range("J6").Value = 0.03972142724517488
$ws.Range("K6").Value = 0.7654856254559372
$ws.Range("L6").Value = 0.7886132731633581
$ws.Range("M6").Value = 0.828811572615664
$ws.Range("N6").Value = 0.8063106808976204
$ws.Range("O6").Value = 0.8819159684619953
$ws.Range("P6").Value = 0.7609529102436011
$ws.Range("Q6").Value = 0.02205912082744699
$ws.Range("R6").Value = 0.748255198087225
$ws.Range("S6").Value = 0.7566536360891201
$ws.Range("T6").Value = 0.7966758337824428
$ws.Range("U6").Value = 0.7315668202764978
$ws.Range("V6").Value = 0.7716130629827197
$ws.Range("W6").Value = 0.8082327911765731
$ws.Range("X6").Value = 0.02678402651696709
$ws.Range("Y6").Value = 0.7966758396235459
$ws.Range("Z6").Value = 0.7796123624908123
$ws.Range("AA6").Value = 0.8301499180156344
$ws.Range("AB6").Value = 0.7858418975966166
$ws.Range("AC6").Value = 0.8488839381562563
$ws.Range("AD6").Value = 0.8079642493184324
$ws.Range("AE6").Value = 0.02998949629178382
$ws.Range("AF6").Value = 0.799912498848669
$ws.Range("AG6").Value = 0.7669725649059193
$ws.Range("AH6").Value = 0.8184277682672393
$ws.Range("AI6").Value = 0.7964482794703507
$ws.Range("AJ6").Value = 0.8580601350999834
$ws.Range("AK6").Value = 0.8169122273468202
$ws.Range("AL6").Value = 0.04888747874274001
$ws.Range("AM6").Value = 0.8057839426923914
$ws.Range("AN6").Value = 0.7322024802100703
$ws.Range("AO6").Value = 0.8276870232515393
$ws.Range("AP6").Value = 0.8380385348127284
$ws.Range("AQ6").Value = 0.8808491557673718
$ws.Range("B7").Value = 0.8526172370702561
$ws.Range("C7").Value = 0.04726931293461654
$ws.Range("D7").Value = 0.8196711441392293
$ws.Range("E7").Value = 0.8084052148568277
$ws.Range("F7").Value = 0.8931043850005298
$ws.Range("G7").Value = 0.8170073653944623
$ws.Range("H7").Value = 0.9248980759602318
$ws.Range("I7").Value = 0.8715642118514711
$ws.Range("J7").Value = 0.04207819240748639
$ws.Range("K7").Value = 0.8509001636661211
$ws.Range("L7").Value = 0.818456812005199
$ws.Range("M7").Value = 0.9036565568430243
$ws.Range("N7").Value = 0.8493515848354559
$ws.Range("P7").Value = 0.8305318399463862
$ws.Range("Q7").Value = 0.04252656265324263
$ws.Range("R7").Value = 0.8106669019913335
$ws.Range("S7").Value = 0.7998279410946603
$ws.Range("U7").Value = 0.8166753185507876
$ws.Range("W7").Value = 0.8795942503981061
$ws.Range("X7").Value = 0.04652763315824815
$ws.Range("Y7").Value = 0.8923304374134681
$ws.Range("AB7").Value = 0.8710465162078066
$ws.Range("AC7").Value = 0.90248223605609
$ws.Range("AK7").Value = 0.86305094578302
$ws.Range("AL7").Value = 0.05070449712042256
$ws.Range("AM7").Value = 0.9046190216402983
$ws.Range("AN7").Value = 0.7760878623357402
$ws.Range("AO7").Value = 0.9146506683464736
$ws.Range("AP7").Value = 0.8389430328946458
$ws.Range("AQ7").Value = 0.8809541436979426
